# Add new columns I (I0) and J (IF) to match the updated report layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style/formatting from the existing "IP" header (H1)
# onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows: (row, I value, J value)
$data = @(
    @(2, 1, 6),
    @(3, 1, 5),
    @(4, 1, 6),
    @(5, 1, 6),
    @(6, 1, 6),
    @(7, 1, 7),
    @(8, 4, 7),
    @(9, 3, 7),
    @(10, 1, 5),
    @(11, 1, 6),
    @(12, 1, 3),
    @(13, 1, 6),
    @(14, 1, 6),
    @(15, 1, 6),
    @(16, 1, 6),
    @(17, 1, 6),
    @(18, 1, 6),
    @(19, 1, 7),
    @(20, 1, 6),
    @(21, 1, 4),
    @(22, 1, 5),
    @(23, 1, 6),
    @(24, 1, 6),
    @(25, 1, 5),
    @(26, 1, 4),
    @(27, 1, 5),
    @(28, 1, 3),
    @(29, 1, 2)
)

foreach ($row in $data) {
    $r = $row[0]
    $i = $row[1]
    $j = $row[2]
    $ws.Cells.Item($r, 9).Value = $i
    $ws.Cells.Item($r, 10).Value = $j
}
